# Updated full name functions
#
# Replace the verbose `xxx.name.full(middle=’full’)` Jinja expressions with
# the simpler `xxx.name_full()` helper in four specific spots of the
# revocation notice paragraph. One occurrence of the pattern (after
# "empowering {{ ") and one `person.name.full(middle = ‘full’)` occurrence
# (with spaces around `=`, "shall remain a successor agent") are left as-is.

$d = $word.ActiveDocument

$lsq = [char]0x2018
$rsq = [char]0x2019

# 1) "I, {{ user.name.full(middle=’full’) }} of {{ ..."
$old1 = "user.name.full(middle=" + $rsq + "full" + $rsq + ")"
$new1 = "user.name_full()"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# 2) "{{ new_property_agent.name.full(middle=’full’) }} shall now take the place of"
$old2 = "new_property_agent.name.full(middle=" + $rsq + "full" + $rsq + ")"
$new2 = "new_property_agent.name_full()"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# 3) "... shall now take the place of {{ property_agent.name.full(middle=’full’) }} as my agent for Power of Attorney for Property."
$old3 = "property_agent.name.full(middle=" + $rsq + "full" + $rsq + ") }} as my agent for Power"
$new3 = "property_agent.name_full() }} as my agent for Power"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# 4) "... granted to {{ person.name.full(middle=’full’) }}."
$old4 = "person.name.full(middle=" + $rsq + "full" + $rsq + ") }}."
$new4 = "person.name_full() }}."
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
